$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O mirrors column N's per-row formatting for rows 3-14, with new 2021 data.
# Copy format row-by-row from N to O so each row keeps its own style (not a single
# uniform style from a block copy).
for ($r = 3; $r -le 14; $r++) {
  $ws.Range("N$r").Copy()
  $ws.Range("O$r").PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Set the new 2021 values for column O
$ws.Cells.Item(4, 15).Value = 2021
$ws.Cells.Item(5, 15).Value = 97
$ws.Cells.Item(6, 15).Value = 96.2
$ws.Cells.Item(7, 15).Value = 62.7
$ws.Cells.Item(8, 15).Value = 100
$ws.Cells.Item(9, 15).Value = 100
$ws.Cells.Item(10, 15).Value = "-"
$ws.Cells.Item(11, 15).Value = 100
$ws.Cells.Item(12, 15).Value = 57.9
$ws.Cells.Item(13, 15).Value = 100
$ws.Cells.Item(14, 15).Value = "-"

# Update selection to match diff (activeCell O17)
$ws.Range("O17").Select()
